$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'time_taken',
    '2021-10-05 13:41:54.684330',
    '2021-10-05 13:41:54.684342',
    '2021-10-05 13:41:54.684345',
    '2021-10-05 13:41:54.684348',
    '2021-10-05 13:41:54.684351',
    '2021-10-05 13:41:54.684353',
    '2021-10-05 13:41:54.684356',
    '2021-10-05 13:41:54.684359',
    '2021-10-05 13:41:54.684362',
    '2021-10-05 13:41:54.684364',
    '2021-10-05 13:41:54.684367',
    '2021-10-05 13:41:54.684370',
    '2021-10-05 13:41:54.684372',
    '2021-10-05 13:41:54.684375',
    '2021-10-05 13:41:54.684377',
    '2021-10-05 13:41:54.684380',
    '2021-10-05 13:41:54.684383',
    '2021-10-05 13:41:54.684385',
    '2021-10-05 13:41:54.684388',
    '2021-10-05 13:41:54.684391',
    '2021-10-05 13:41:54.684393',
    '2021-10-05 13:41:54.684396',
    '2021-10-05 13:41:54.684398',
    '2021-10-05 13:41:54.684401',
    '2021-10-05 13:41:54.684404',
    '2021-10-05 13:41:54.684407',
    '2021-10-05 13:41:54.684409',
    '2021-10-05 13:41:54.684412',
    '2021-10-05 13:41:54.684414',
    '2021-10-05 13:41:54.684417',
    '2021-10-05 13:41:54.684419',
    '2021-10-05 13:41:54.684422',
    '2021-10-05 13:41:54.684425',
    '2021-10-05 13:41:54.684428',
    '2021-10-05 13:41:54.684430',
    '2021-10-05 13:41:54.684433',
    '2021-10-05 13:41:54.684435',
    '2021-10-05 13:41:54.684438',
    '2021-10-05 13:41:54.684441',
    '2021-10-05 13:41:54.684443',
    '2021-10-05 13:41:54.684446',
    '2021-10-05 13:41:54.684449',
    '2021-10-05 13:41:54.684451',
    '2021-10-05 13:41:54.684454',
    '2021-10-05 13:41:54.684456',
    '2021-10-05 13:41:54.684459',
    '2021-10-05 13:41:54.684461',
    '2021-10-05 13:41:54.684464',
    '2021-10-05 13:41:54.684466',
    '2021-10-05 13:41:54.684469',
    '2021-10-05 13:41:54.684472',
    '2021-10-05 13:41:54.684474',
    '2021-10-05 13:41:54.684477',
    '2021-10-05 13:41:54.684480',
    '2021-10-05 13:41:54.684483',
    '2021-10-05 13:41:54.684485',
    '2021-10-05 13:41:54.684488',
    '2021-10-05 13:41:54.684490',
    '2021-10-05 13:41:54.684493',
    '2021-10-05 13:41:54.684496',
    '2021-10-05 13:41:54.684498',
    '2021-10-05 13:41:54.684501',
    '2021-10-05 13:41:54.684503',
    '2021-10-05 13:41:54.684506',
    '2021-10-05 13:41:54.684510',
    '2021-10-05 13:41:54.684513',
    '2021-10-05 13:41:54.684515',
    '2021-10-05 13:41:54.684518',
    '2021-10-05 13:41:54.684521',
    '2021-10-05 13:41:54.684523',
    '2021-10-05 13:41:54.684526',
    '2021-10-05 13:41:54.684529',
    '2021-10-05 13:41:54.684531',
    '2021-10-05 13:41:54.684534',
    '2021-10-05 13:41:54.684536',
    '2021-10-05 13:41:54.684539',
    '2021-10-05 13:41:54.684544',
    '2021-10-05 13:41:54.684547',
    '2021-10-05 13:41:54.684550',
    '2021-10-05 13:41:54.684553',
    '2021-10-05 13:41:54.684555',
    '2021-10-05 13:41:54.684558',
    '2021-10-05 13:41:54.684560',
    '2021-10-05 13:41:54.684563',
    '2021-10-05 13:41:54.684565',
    '2021-10-05 13:41:54.684568',
    '2021-10-05 13:41:54.684570',
    '2021-10-05 13:41:54.684573',
    '2021-10-05 13:41:54.684575',
    '2021-10-05 13:41:54.684578',
    '2021-10-05 13:41:54.684580',
    '2021-10-05 13:41:54.684583',
    '2021-10-05 13:41:54.684587',
    '2021-10-05 13:41:54.684590',
    '2021-10-05 13:41:54.684592',
    '2021-10-05 13:41:54.684595',
    '2021-10-05 13:41:54.684597',
    '2021-10-05 13:41:54.684600',
    '2021-10-05 13:41:54.684602',
    '2021-10-05 13:41:54.684605',
    '2021-10-05 13:41:54.684608',
    '2021-10-05 13:41:54.684610',
    '2021-10-05 13:41:54.684613',
    '2021-10-05 13:41:54.684615',
    '2021-10-05 13:41:54.684618',
    '2021-10-05 13:41:54.684620',
    '2021-10-05 13:41:54.684623',
    '2021-10-05 13:41:54.684626',
    '2021-10-05 13:41:54.684630',
    '2021-10-05 13:41:54.684633',
    '2021-10-05 13:41:54.684636',
    '2021-10-05 13:41:54.684638',
    '2021-10-05 13:41:54.684641',
    '2021-10-05 13:41:54.684643',
    '2021-10-05 13:41:54.684646',
    '2021-10-05 13:41:54.684648',
    '2021-10-05 13:41:54.684651',
    '2021-10-05 13:41:54.684653',
    '2021-10-05 13:41:54.684656',
    '2021-10-05 13:41:54.684658',
    '2021-10-05 13:41:54.684661',
    '2021-10-05 13:41:54.684664',
    '2021-10-05 13:41:54.684666',
    '2021-10-05 13:41:54.684669',
    '2021-10-05 13:41:54.684671',
    '2021-10-05 13:41:54.684674',
    '2021-10-05 13:41:54.684676',
    '2021-10-05 13:41:54.684679',
    '2021-10-05 13:41:54.684683',
    '2021-10-05 13:41:54.684686',
    '2021-10-05 13:41:54.684689',
    '2021-10-05 13:41:54.684692',
    '2021-10-05 13:41:54.684694',
    '2021-10-05 13:41:54.684697',
    '2021-10-05 13:41:54.684699',
    '2021-10-05 13:41:54.684702',
    '2021-10-05 13:41:54.684704',
    '2021-10-05 13:41:54.684707',
    '2021-10-05 13:41:54.684710',
    '2021-10-05 13:41:54.684712',
    '2021-10-05 13:41:54.684715',
    '2021-10-05 13:41:54.684717',
    '2021-10-05 13:41:54.684720',
    '2021-10-05 13:41:54.684722',
    '2021-10-05 13:41:54.684725',
    '2021-10-05 13:41:54.684728',
    '2021-10-05 13:41:54.684730',
    '2021-10-05 13:41:54.684733',
    '2021-10-05 13:41:54.684736',
    '2021-10-05 13:41:54.684739',
    '2021-10-05 13:41:54.684741',
    '2021-10-05 13:41:54.684744',
    '2021-10-05 13:41:54.684746',
    '2021-10-05 13:41:54.684749',
    '2021-10-05 13:41:54.684752',
    '2021-10-05 13:41:54.684754',
    '2021-10-05 13:41:54.684757',
    '2021-10-05 13:41:54.684760',
    '2021-10-05 13:41:54.684762',
    '2021-10-05 13:41:54.684765',
    '2021-10-05 13:41:54.684767',
    '2021-10-05 13:41:54.684770',
    '2021-10-05 13:41:54.684773',
    '2021-10-05 13:41:54.684775',
    '2021-10-05 13:41:54.684778',
    '2021-10-05 13:41:54.684780',
    '2021-10-05 13:41:54.684783',
    '2021-10-05 13:41:54.684785',
    '2021-10-05 13:41:54.684788',
    '2021-10-05 13:41:54.684791',
    '2021-10-05 13:41:54.684793',
    '2021-10-05 13:41:54.684796',
    '2021-10-05 13:41:54.684800',
    '2021-10-05 13:41:54.684803',
    '2021-10-05 13:41:54.684806',
    '2021-10-05 13:41:54.684809',
    '2021-10-05 13:41:54.684812',
    '2021-10-05 13:41:54.684814',
    '2021-10-05 13:41:54.684817',
    '2021-10-05 13:41:54.684820',
    '2021-10-05 13:41:54.684823',
    '2021-10-05 13:41:54.684825',
    '2021-10-05 13:41:54.684828',
    '2021-10-05 13:41:54.684831',
    '2021-10-05 13:41:54.684833',
    '2021-10-05 13:41:54.684836',
    '2021-10-05 13:41:54.684839',
    '2021-10-05 13:41:54.684842',
    '2021-10-05 13:41:54.684844',
    '2021-10-05 13:41:54.684847',
    '2021-10-05 13:41:54.684849',
    '2021-10-05 13:41:54.684852',
    '2021-10-05 13:41:54.684855',
    '2021-10-05 13:41:54.684857',
    '2021-10-05 13:41:54.684860',
    '2021-10-05 13:41:54.684862',
    '2021-10-05 13:41:54.684865',
    '2021-10-05 13:41:54.684868',
    '2021-10-05 13:41:54.684870',
    '2021-10-05 13:41:54.684873',
    '2021-10-05 13:41:54.684875',
    '2021-10-05 13:41:54.684878',
    '2021-10-05 13:41:54.684881',
    '2021-10-05 13:41:54.684883',
    '2021-10-05 13:41:54.684886',
    '2021-10-05 13:41:54.684888',
    '2021-10-05 13:41:54.684891',
    '2021-10-05 13:41:54.684894',
    '2021-10-05 13:41:54.684896',
    '2021-10-05 13:41:54.684899',
    '2021-10-05 13:41:54.684901',
    '2021-10-05 13:41:54.684904',
    '2021-10-05 13:41:54.684907',
    '2021-10-05 13:41:54.684909',
    '2021-10-05 13:41:54.684912',
    '2021-10-05 13:41:54.684915',
    '2021-10-05 13:41:54.684917',
    '2021-10-05 13:41:54.684920',
    '2021-10-05 13:41:54.684922',
    '2021-10-05 13:41:54.684925',
    '2021-10-05 13:41:54.684927',
    '2021-10-05 13:41:54.684930',
    '2021-10-05 13:41:54.684933',
    '2021-10-05 13:41:54.684935',
    '2021-10-05 13:41:54.684938',
    '2021-10-05 13:41:54.684941',
    '2021-10-05 13:41:54.684943',
    '2021-10-05 13:41:54.684946',
    '2021-10-05 13:41:54.684948',
    '2021-10-05 13:41:54.684951',
    '2021-10-05 13:41:54.684954',
    '2021-10-05 13:41:54.684956',
    '2021-10-05 13:41:54.684960',
    '2021-10-05 13:41:54.684963',
    '2021-10-05 13:41:54.684966',
    '2021-10-05 13:41:54.684969',
    '2021-10-05 13:41:54.684971',
    '2021-10-05 13:41:54.684974',
    '2021-10-05 13:41:54.684976',
    '2021-10-05 13:41:54.684979',
    '2021-10-05 13:41:54.684981',
    '2021-10-05 13:41:54.684984',
    '2021-10-05 13:41:54.684987',
    '2021-10-05 13:41:54.684989',
    '2021-10-05 13:41:54.684992',
    '2021-10-05 13:41:54.684994',
    '2021-10-05 13:41:54.684997',
    '2021-10-05 13:41:54.685000',
    '2021-10-05 13:41:54.685002',
    '2021-10-05 13:41:54.685005',
    '2021-10-05 13:41:54.685008',
    '2021-10-05 13:41:54.685010',
    '2021-10-05 13:41:54.685013',
    '2021-10-05 13:41:54.685015',
    '2021-10-05 13:41:54.685018',
    '2021-10-05 13:41:54.685020',
    '2021-10-05 13:41:54.685023',
    '2021-10-05 13:41:54.685025',
    '2021-10-05 13:41:54.685028',
    '2021-10-05 13:41:54.685031',
    '2021-10-05 13:41:54.685033',
    '2021-10-05 13:41:54.685036',
    '2021-10-05 13:41:54.685039',
    '2021-10-05 13:41:54.685041',
    '2021-10-05 13:41:54.685044',
    '2021-10-05 13:41:54.685046',
    '2021-10-05 13:41:54.685049',
    '2021-10-05 13:41:54.685051',
    '2021-10-05 13:41:54.685054',
    '2021-10-05 13:41:54.685057',
    '2021-10-05 13:41:54.685059',
    '2021-10-05 13:41:54.685062',
    '2021-10-05 13:41:54.685065',
    '2021-10-05 13:41:54.685067',
    '2021-10-05 13:41:54.685070',
    '2021-10-05 13:41:54.685072',
    '2021-10-05 13:41:54.685075',
    '2021-10-05 13:41:54.685077',
    '2021-10-05 13:41:54.685080',
    '2021-10-05 13:41:54.685082',
    '2021-10-05 13:41:54.685085',
    '2021-10-05 13:41:54.685088',
    '2021-10-05 13:41:54.685090',
    '2021-10-05 13:41:54.685093',
    '2021-10-05 13:41:54.685095',
    '2021-10-05 13:41:54.685098',
    '2021-10-05 13:41:54.685101',
    '2021-10-05 13:41:54.685103',
    '2021-10-05 13:41:54.685106',
    '2021-10-05 13:41:54.685108',
    '2021-10-05 13:41:54.685111',
    '2021-10-05 13:41:54.685114',
    '2021-10-05 13:41:54.685116',
    '2021-10-05 13:41:54.685119',
    '2021-10-05 13:41:54.685121',
    '2021-10-05 13:41:54.685124',
    '2021-10-05 13:41:54.685127',
    '2021-10-05 13:41:54.685129',
    '2021-10-05 13:41:54.685132',
    '2021-10-05 13:41:54.685134',
    '2021-10-05 13:41:54.685137',
    '2021-10-05 13:41:54.685140',
    '2021-10-05 13:41:54.685142',
    '2021-10-05 13:41:54.685145',
    '2021-10-05 13:41:54.685147',
    '2021-10-05 13:41:54.685150',
    '2021-10-05 13:41:54.685153',
    '2021-10-05 13:41:54.685155',
    '2021-10-05 13:41:54.685159',
    '2021-10-05 13:41:54.685162',
    '2021-10-05 13:41:54.685165',
    '2021-10-05 13:41:54.685177',
    '2021-10-05 13:41:54.685180',
    '2021-10-05 13:41:54.685183',
    '2021-10-05 13:41:54.685185',
    '2021-10-05 13:41:54.685188',
    '2021-10-05 13:41:54.685190',
    '2021-10-05 13:41:54.685193',
    '2021-10-05 13:41:54.685196',
    '2021-10-05 13:41:54.685199',
    '2021-10-05 13:41:54.685201',
    '2021-10-05 13:41:54.685204',
    '2021-10-05 13:41:54.685207'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 6).Value = $values[$i]
}

$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

Write-Host "Added time_taken column"